$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CAPEX")
$ws.Range("U6").Value = 41616.833333333336
